$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stack")

# Fix the TC-Remove column (I) time complexity values: O(1) -> O(n)
# Data rows are 2 through 101 for column I.
for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 9).Value = "O(n)"
}

# Update the view state: scroll to row 80 and select I2:I101 with active cell I2
$ws.Range("I2:I101").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 80
